$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells, copying the existing header style (H1) so
# the new columns match the other headers (bold, bordered, centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 / IF values for rows 2-41 (A2:A41 => dataset rows 0-39)
$data = @(
    @(1, 5),
    @(1, 6),
    @(4, 7),
    @(3, 8),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 4),
    @(2, 6),
    @(1, 5),
    @(1, 4),
    @(7, 8),
    @(4, 8),
    @(6, 8),
    @(6, 8),
    @(1, 2),
    @(11, 12),
    @(7, 9),
    @(7, 9),
    @(6, 7),
    @(7, 7),
    @(1, 1),
    @(2, 4),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 7),
    @(1, 3),
    @(1, 6),
    @(1, 4),
    @(1, 5),
    @(1, 4),
    @(1, 2),
    @(5, 5)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}
